# Update the binning information table to reflect a log-scale dependent
# variable: fewer, wider bins for TENURE/HRLYEARN/UHRSMAIN.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..14 (Feature, Bin Index, Bin Start, Bin End)
$data = @(
    @("UHRSMAIN", 0, 33.8,   37.5),
    @("UHRSMAIN", 1, 37.5,   40),
    @("UHRSMAIN", 2, 40,     43.5),
    @("TENURE",   0, 1,      20),
    @("TENURE",   1, 20,     52),
    @("TENURE",   2, 52,     105),
    @("TENURE",   3, 105,    199),
    @("TENURE",   4, 199,    240),
    @("HRLYEARN", 0, 5.77,   22.83400000000001),
    @("HRLYEARN", 1, 22.83400000000001, 28.85),
    @("HRLYEARN", 2, 28.85,  37.02),
    @("HRLYEARN", 3, 37.02,  48.08),
    @("HRLYEARN", 4, 48.08,  79.37)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

# The table shrank from 16 data+header rows to 14; clear the now-unused
# rows 15 and 16 so the sheet's used range collapses back to A1:D14.
$ws.Range("A15:D16").Clear()
